$wb = $excel.ActiveWorkbook

# Sheet 1: weibull
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -1.95308165608033
$ws.Range("C2").Value = 0.235764186065949
$ws.Range("B3").Value = -0.127807148941032
$ws.Range("C3").Value = 0.128157101154197

# Sheet 2: lognormal
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 1.05257362849532
$ws.Range("C2").Value = 0.288080901041177
$ws.Range("B3").Value = -0.717436925236619
$ws.Range("C3").Value = 0.124537607141585

# Sheet 3: llogis
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.50596150897498
$ws.Range("C2").Value = 0.249752256684582
$ws.Range("B3").Value = 0.222267720042517
$ws.Range("C3").Value = 0.166463011139801

# Sheet 4: gompertz
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -1.8334625702322
$ws.Range("C2").Value = 0.270318049401206
$ws.Range("B3").Value = -0.0444484446995318
$ws.Range("C3").Value = 0.0368542173565487

# Sheet 6: weibull cov
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0555847514313396
$ws.Range("B2").Value = -0.014699123668472
$ws.Range("A3").Value = -0.014699123668472
$ws.Range("B3").Value = 0.016424242576247

# Sheet 7: lognormal cov
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0829906055446962
$ws.Range("B2").Value = -0.0278016389585685
$ws.Range("A3").Value = -0.0278016389585685
$ws.Range("B3").Value = 0.0155096155925518

# Sheet 8: llogis cov
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0623761897190414
$ws.Range("B2").Value = -0.00515070683471067
$ws.Range("A3").Value = -0.00515070683471067
$ws.Range("B3").Value = 0.0277099340777296

# Sheet 9: gompertz cov
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0730718478320728
$ws.Range("B2").Value = -0.00680553207578169
$ws.Range("A3").Value = -0.00680553207578169
$ws.Range("B3").Value = 0.00135823333696374
